$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Initial Estimate" values for the first block of tasks (rows 4-12,
# column D) are replaced with a literal "/" text marker instead of a
# numeric estimate (and two previously-blank cells, D8/D9, also get the
# same marker).
$slashRange = $ws.Range("D4:D12")
$slashRange.Value = "/"

# Recalculate so the dependent SUM formulas (D26:G26) and the chart that
# is driven off them pick up the new totals.
$excel.CalculateFull()

# Refresh the burndown chart so its cached values follow the recalculated
# totals in C26:G26.
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$chart.Refresh()

# The active selection moved from E8 to B15 before the file was saved.
$ws.Range("B15").Select()
